$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cadastroUsuario")
$ws.Range("A2").Value = "tesTe08"
